$d = $word.ActiveDocument

$replacements = @(
    @{Old = "2023-10-11 Wednesday"; New = "2023-10-12 Thursday"},
    @{Old = "14×53=742";  New = "82×61=5002"},
    @{Old = "97×56=5432"; New = "16×91=1456"},
    @{Old = "91×24=2184"; New = "49×40=1960"},
    @{Old = "57×59=3363"; New = "53×39=2067"},
    @{Old = "30×30=900";  New = "62×38=2356"},
    @{Old = "64×50=3200"; New = "81×19=1539"},
    @{Old = "35×43=1505"; New = "49×38=1862"},
    @{Old = "96×72=6912"; New = "12×45=540"},
    @{Old = "75×42=3150"; New = "81×58=4698"},
    @{Old = "85×78=6630"; New = "25×79=1975"},
    @{Old = "85×16=1360"; New = "59×37=2183"},
    @{Old = "71×69=4899"; New = "83×88=7304"},
    @{Old = "48×49=2352"; New = "30×87=2610"},
    @{Old = "54×83=4482"; New = "25×35=875"},
    @{Old = "47×71=3337"; New = "92×98=9016"},
    @{Old = "60×50=3000"; New = "76×22=1672"},
    @{Old = "55×89=4895"; New = "78×46=3588"},
    @{Old = "69×19=1311"; New = "92×44=4048"},
    @{Old = "14×67=938";  New = "85×66=5610"},
    @{Old = "39×63=2457"; New = "98×77=7546"},
    @{Old = "80×87=6960"; New = "20×18=360"},
    @{Old = "34×78=2652"; New = "62×71=4402"},
    @{Old = "51×54=2754"; New = "40×75=3000"},
    @{Old = "36×18=648";  New = "22×93=2046"},
    @{Old = "69×79=5451"; New = "69×52=3588"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2)
}
